$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "1`nEG-202-01807-McInnis-IDE-318`nEG-202-01806-Mohamed Ali-IDE-206A`nEG-202-01806-Mohamed Ali-IDE-206B"
$ws.Range("E2").Value = "13"
$ws.Range("B3").Value = "3`nEG-350-03807-McInnis-IDE-206A`nEG-350-03807-McInnis-IDE-206B`nEG-318-03805-Guo-IDE-323`nEG-208-03804-nan-IDE-214`nEG-208-03804-nan-IDE-217A"
$ws.Range("C3").Value = "4`nEG-310-04813-Guo-IDE-113A`nEG-209-04812-Beltramo-IDE-107A`nEG-209-04812-Beltramo-IDE-107B"
$ws.Range("D3").Value = "3`nEG-350-03807-McInnis-IDE-206A`nEG-350-03807-McInnis-IDE-206B`nEG-318-03805-Guo-IDE-323`nEG-208-03804-nan-IDE-214`nEG-208-03804-nan-IDE-217A"
$ws.Range("E3").Value = "14`nEG-498-14808-Rosner-IDE-113A`nEG-498-14807-Moghimi-IDE-128B`nEG-498-14806-Daigneau-IDE-128A"
$ws.Range("F3").Value = "4.0`nEG-310-04813-Guo-IDE-113A`nEG-209-04812-Beltramo-IDE-107A`nEG-209-04812-Beltramo-IDE-107B"
$ws.Range("B4").Value = "5`nEG-412-05809-Sadraey-IDE-317`nEE-318-05808-Moghimi-IDE-206A`nEE-318-05808-Moghimi-IDE-206B`nEG-325-05814-Gorczyca-IDE-107A`nEG-325-05814-Gorczyca-IDE-107B`nEG-207-05811-Daigneau-IDE-217A"
$ws.Range("C4").Value = "13`nEG-419FT-13809-Rosner-IDE-318`nEE-318L-13806-Moghimi-IDE-206A`nEE-318L-13806-Moghimi-IDE-206B"
$ws.Range("D4").Value = "6`nEG-424-06807-Guo-EAX-019`nEG-424-06807-Guo-IDE-107A`nEG-424-06807-Guo-IDE-107B`nEG-360-06811-Sadraey-IDE-204A`nEG-361-06810-Sadique-IDE-118A`nEG-361-06810-Sadique-IDE-118B"
$ws.Range("E4").Value = "5`nEG-412-05809-Sadraey-IDE-317`nEE-318-05808-Moghimi-IDE-206A`nEE-318-05808-Moghimi-IDE-206B`nEG-325-05814-Gorczyca-IDE-107A`nEG-325-05814-Gorczyca-IDE-107B`nEG-207-05811-Daigneau-IDE-217A"
$ws.Range("F4").Value = "6.0`nEG-424-06807-Guo-EAX-019`nEG-424-06807-Guo-IDE-107A`nEG-424-06807-Guo-IDE-107B`nEG-360-06811-Sadraey-IDE-204A`nEG-361-06810-Sadique-IDE-118A`nEG-361-06810-Sadique-IDE-118B"
$ws.Range("B5").Value = "7`nEG-314-07807-Guo-IDE-118A`nEG-314-07807-Guo-IDE-118B`nEG-340-07809-nan-IDE-107A`nEG-340-07809-nan-IDE-107B"
$ws.Range("C5").Value = "14`nEG-498-14808-Rosner-IDE-113A`nEG-498-14807-Moghimi-IDE-128B`nEG-498-14806-Daigneau-IDE-128A"
$ws.Range("D5").Value = "8`nEG-419-08810-Rosner-IDE-318`nEG-308-08813-Sadraey-IDE-323`nEG-316-08809-Moghimi-IDE-107A`nEG-316-08809-Moghimi-IDE-107B`nEG-201-08814-nan-IDE-118A`nEG-201-08814-nan-IDE-118B`nEG-110-08803-nan-nan-nan`nEG-110-08811-nan-nan-nan"
$ws.Range("E5").Value = "7`nEG-314-07807-Guo-IDE-118A`nEG-314-07807-Guo-IDE-118B`nEG-340-07809-nan-IDE-107A`nEG-340-07809-nan-IDE-107B"
$ws.Range("F5").Value = "8.0`nEG-419-08810-Rosner-IDE-318`nEG-308-08813-Sadraey-IDE-323`nEG-316-08809-Moghimi-IDE-107A`nEG-316-08809-Moghimi-IDE-107B`nEG-201-08814-nan-IDE-118A`nEG-201-08814-nan-IDE-118B`nEG-110-08803-nan-nan-nan`nEG-110-08811-nan-nan-nan"
$ws.Range("B6").Value = "9`nEE-310-09808-Moghimi-IDE-206A`nEE-310-09808-Moghimi-IDE-206B`nEG-200-09811-Gorczyca-IDE-118A`nEG-200-09811-Gorczyca-IDE-118B`nEG-207-09813-Daigneau-IDE-217A"
$ws.Range("C6").Value = "10`nEG-110-10810-Kolenbrander-IDE-128A`nEG-110-10810-Kolenbrander-IDE-128B`nEG-110-10809-Eshed-IDE-128A`nEG-110-10809-Eshed-IDE-128B`nEG-335-10811-Daigneau-IDE-206A`nEG-335-10811-Daigneau-IDE-206B`nEG-409-10812-nan-nan-nan"
$ws.Range("E6").Value = "9`nEE-310-09808-Moghimi-IDE-206A`nEE-310-09808-Moghimi-IDE-206B`nEG-200-09811-Gorczyca-IDE-118A`nEG-200-09811-Gorczyca-IDE-118B`nEG-207-09813-Daigneau-IDE-217A"
$ws.Range("F6").Value = "10.0`nEG-110-10810-Kolenbrander-IDE-128A`nEG-110-10810-Kolenbrander-IDE-128B`nEG-110-10809-Eshed-IDE-128A`nEG-110-10809-Eshed-IDE-128B`nEG-335-10811-Daigneau-IDE-206A`nEG-335-10811-Daigneau-IDE-206B`nEG-409-10812-nan-nan-nan"
$ws.Range("B7").Value = "11`nEG-203-11813-McInnis-IDE-323`nEG-200-11806-Gorczyca-IDE-118A`nEG-200-11806-Gorczyca-IDE-118B"
$ws.Range("D7").Value = "12`nEG-410-12808-Sadique-IDE-318`nEG-316-12806-nan-IDE-118A`nEG-316-12806-nan-IDE-118B"
$ws.Range("E7").Value = "11`nEG-203-11813-McInnis-IDE-323`nEG-200-11806-Gorczyca-IDE-118A`nEG-200-11806-Gorczyca-IDE-118B"
$ws.Range("F7").Value = "12.0`nEG-410-12808-Sadique-IDE-318`nEG-316-12806-nan-IDE-118A`nEG-316-12806-nan-IDE-118B"
$ws.Range("B8").Value = "20`nEE-305-20803-Benzerrouk-IDE-206A`nEE-305-20803-Benzerrouk-IDE-206B"
$ws.Range("C8").Value = "21`nEG-110-21805-nan-IDE-217A`nEG-207-21806-nan-IDE-107B"
$ws.Range("D8").Value = "20`nEE-305-20803-Benzerrouk-IDE-206A`nEE-305-20803-Benzerrouk-IDE-206B"
$ws.Range("E8").Value = "21`nEG-110-21805-nan-IDE-217A`nEG-207-21806-nan-IDE-107B"
